$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "K" column (G) had been populated from a stale "Strike#" field.
# Regenerate it with the corrected strikeout counts per start, then
# recompute the downstream std/mean summary and rewrite the s_vals.
$kByRow = @{
    2 = 0
    3 = 0
    4 = 2
    5 = 0
    6 = 0
    7 = 2
    8 = 2
    9 = 1
    10 = 0
    11 = 0
    12 = 0
    13 = 1
    14 = 1
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 1
    20 = 0
    21 = 1
    22 = 3
    23 = 1
    24 = 0
    25 = 0
    27 = 2
    28 = 0
    29 = 2
    30 = 1
    31 = 1
    32 = 2
    33 = 1
    34 = 0
    35 = 1
    36 = 0
    37 = 0
    38 = 1
    39 = 1
    41 = 3
    42 = 2
    43 = 0
    44 = 1
    45 = 1
    46 = 0
    47 = 2
    48 = 2
    49 = 3
    50 = 0
    51 = 1
    52 = 3
    53 = 2
    54 = 4
    55 = 1
    56 = 1
    57 = 3
    58 = 2
    59 = 3
    60 = 0
    61 = 2
    62 = 0
    63 = 3
    64 = 2
    65 = 1
    66 = 1
    67 = 3
    69 = 2
    70 = 3
    71 = 1
    72 = 2
    73 = 3
    74 = 1
    75 = 3
    76 = 2
    77 = 2
}

foreach ($row in $kByRow.Keys) {
    $ws.Cells.Item($row, 7).Value = $kByRow[$row]
}
